$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base date in B2; formulas in B3:B21 reference B2 (directly or
# indirectly) and will recalculate automatically.
$ws.Range("B2").Value = 43556

# Update the active selection to match the author's saved view state.
$ws.Range("B3").Select()
